# Auto-generated edit script applying scraped market-data refresh to Sophia_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2700.926
$ws.Range("I15").Value = 2700.926
$ws.Range("K15").Value = 8102.778
$ws.Range("M15").Value = -7933.778

$ws.Range("H51").Value = 1444.4445
$ws.Range("I51").Value = 1000
$ws.Range("K51").Value = 1000
$ws.Range("M51").Value = -516

$ws.Range("H58").Value = 1452.8572
$ws.Range("I58").Value = 1364.3334
$ws.Range("K58").Value = 4093.0002
$ws.Range("M58").Value = -3943.0002

$ws.Range("H64").Value = 7333.3335
$ws.Range("I64").Value = 7333.3335
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 7333.3335
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -7085.3335
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 7333.3335
$ws.Range("I67").Value = 7333.3335
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 7333.3335
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -6475.3335
$ws.Range("N67").ClearContents()

$ws.Range("H74").Value = 6331.6665
$ws.Range("I74").Value = 6331.6665
$ws.Range("K74").Value = 6331.6665
$ws.Range("M74").Value = -5395.6665

$ws.Range("H76").Value = 1774.25
$ws.Range("I76").Value = 1365.6666
$ws.Range("K76").Value = 1365.6666
$ws.Range("M76").Value = -1050.6666

$ws.Range("H77").Value = 6331.6665
$ws.Range("I77").Value = 6331.6665
$ws.Range("K77").Value = 31658.3325
$ws.Range("M77").Value = -26978.3325

$ws.Range("H79").Value = 1774.25
$ws.Range("I79").Value = 1365.6666
$ws.Range("K79").Value = 1365.6666
$ws.Range("M79").Value = -273.6666

$ws.Range("H82").Value = 5765.6665
$ws.Range("I82").Value = 4625
$ws.Range("J82").Value = 8047
$ws.Range("K82").Value = 13875
$ws.Range("L82").Value = 24141
$ws.Range("M82").Value = -13469
$ws.Range("N82").Value = -24953

$ws.Range("H85").Value = 5765.6665
$ws.Range("I85").Value = 4625
$ws.Range("J85").Value = 8047
$ws.Range("K85").Value = 13875
$ws.Range("L85").Value = 24141
$ws.Range("M85").Value = -12471
$ws.Range("N85").Value = -26949

$ws.Range("H135").Value = 1870.5714
$ws.Range("I135").Value = 619
$ws.Range("K135").Value = 5571
$ws.Range("M135").Value = -3036

$ws.Range("H138").Value = 3757.4
$ws.Range("J138").Value = 4660.4546
$ws.Range("L138").Value = 13981.3638
$ws.Range("N138").Value = -24261.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14269.542
$ws.Range("I32").Value = 14269.542
$ws.Range("K32").Value = 14269.542
$ws.Range("M32").Value = -13982.542

$ws.Range("H61").Value = 1429.4286
$ws.Range("I61").Value = 1641.2
$ws.Range("J61").Value = 900
$ws.Range("K61").Value = 1641.2
$ws.Range("L61").Value = 900
$ws.Range("M61").Value = -1429.2
$ws.Range("N61").Value = -1324

$ws.Range("H63").Value = 10167.5
$ws.Range("J63").Value = 15001.25
$ws.Range("L63").Value = 15001.25
$ws.Range("N63").Value = -16373.25

$ws.Range("H66").Value = 10167.5
$ws.Range("J66").Value = 15001.25
$ws.Range("L66").Value = 75006.25
$ws.Range("N66").Value = -81870.25

$ws.Range("H74").Value = 20363
$ws.Range("I74").Value = 19814.04
$ws.Range("K74").Value = 19814.04
$ws.Range("M74").Value = -18940.04

$ws.Range("H77").Value = 20363
$ws.Range("I77").Value = 19814.04
$ws.Range("K77").Value = 99070.20000000001
$ws.Range("M77").Value = -94702.20000000001

$ws.Range("H132").Value = 2184.913
$ws.Range("I132").Value = 1732.7059
$ws.Range("K132").Value = 5198.1177
$ws.Range("M132").Value = -2668.1177

$ws.Range("H136").Value = 1429.4286
$ws.Range("I136").Value = 1641.2
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 4923.6
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -2373.6
$ws.Range("N136").Value = -7800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6250
$ws.Range("I105").Value = 6250
$ws.Range("K105").Value = 6250
$ws.Range("M105").Value = -4503

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2029.5
$ws.Range("I31").Value = 1659.4
$ws.Range("K31").Value = 1659.4
$ws.Range("M31").Value = -1364.4

$ws.Range("H34").Value = 2029.5
$ws.Range("I34").Value = 1659.4
$ws.Range("K34").Value = 1659.4
$ws.Range("M34").Value = -1457.4

$ws.Range("H51").Value = 50489.5
$ws.Range("J51").Value = 50489.5
$ws.Range("L51").Value = 50489.5
$ws.Range("N51").Value = -51961.5

$ws.Range("H60").Value = 23245.4
$ws.Range("I60").Value = 5064.3335
$ws.Range("J60").Value = 50517
$ws.Range("K60").Value = 5064.3335
$ws.Range("L60").Value = 50517
$ws.Range("M60").Value = -4553.3335
$ws.Range("N60").Value = -51539

$ws.Range("H61").Value = 50489.5
$ws.Range("J61").Value = 50489.5
$ws.Range("L61").Value = 50489.5
$ws.Range("N61").Value = -51185.5

$ws.Range("H132").Value = 1541.3158
$ws.Range("I132").Value = 985.3570999999999
$ws.Range("J132").Value = 3098
$ws.Range("K132").Value = 2956.0713
$ws.Range("L132").Value = 9294
$ws.Range("M132").Value = -426.0712999999996
$ws.Range("N132").Value = -14354

$ws.Range("H134").Value = 4175.5
$ws.Range("I134").Value = 3805.6924
$ws.Range("K134").Value = 11417.0772
$ws.Range("M134").Value = -8882.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1154.2
$ws.Range("I5").Value = 1341.75
$ws.Range("K5").Value = 4025.25
$ws.Range("M5").Value = -3913.25

$ws.Range("H107").Value = 545
$ws.Range("J107").Value = 559.5
$ws.Range("L107").Value = 1678.5
$ws.Range("N107").Value = -5518.5

$ws.Range("H132").Value = 10638
$ws.Range("J132").Value = 24000
$ws.Range("L132").Value = 216000
$ws.Range("N132").Value = -221060

$ws.Range("H135").Value = 1154.2
$ws.Range("I135").Value = 1341.75
$ws.Range("K135").Value = 12075.75
$ws.Range("M135").Value = -9540.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 16667316
$ws.Range("I12").Value = 16667316
$ws.Range("K12").Value = 16667316
$ws.Range("M12").Value = -16667176

$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 2000
$ws.Range("M70").Value = -1730

$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 2000
$ws.Range("M73").Value = -1064

$ws.Range("H122").Value = 7500
$ws.Range("I122").Value = 7500
$ws.Range("K122").Value = 22500
$ws.Range("M122").Value = -20050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5788.0557
$ws.Range("I22").Value = 3957.7273
$ws.Range("J22").Value = 8664.286
$ws.Range("K22").Value = 3957.7273
$ws.Range("L22").Value = 8664.286
$ws.Range("M22").Value = -3662.7273
$ws.Range("N22").Value = -9254.286

$ws.Range("H27").Value = 5788.0557
$ws.Range("I27").Value = 3957.7273
$ws.Range("J27").Value = 8664.286
$ws.Range("K27").Value = 3957.7273
$ws.Range("L27").Value = 8664.286
$ws.Range("M27").Value = -3850.7273
$ws.Range("N27").Value = -8878.286

$ws.Range("H132").Value = 3555.6667
$ws.Range("I132").Value = 3390.111
$ws.Range("K132").Value = 10170.333
$ws.Range("M132").Value = -7640.332999999999

$ws.Range("H136").Value = 3513.7058
$ws.Range("I136").Value = 3430.5715
$ws.Range("K136").Value = 10291.7145
$ws.Range("M136").Value = -7741.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 54950
$ws.Range("J33").Value = 54950
$ws.Range("L33").Value = 54950
$ws.Range("N33").Value = -55450

$ws.Range("H36").Value = 54950
$ws.Range("J36").Value = 54950
$ws.Range("L36").Value = 54950
$ws.Range("N36").Value = -55450

$ws.Range("H132").Value = 1571.1666
$ws.Range("I132").Value = 829.8333
$ws.Range("K132").Value = 2489.4999
$ws.Range("M132").Value = 40.5001000000002

$ws.Range("H136").Value = 1610
$ws.Range("I136").Value = 1610
$ws.Range("K136").Value = 4830
$ws.Range("M136").Value = -2280
